$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook has three 12-row blocks of identical data rows (2-13, 14-25, 26-37).
# This reverts a prior edit that had grown each block from 10 to 12 rows, so we
# shrink each block back down to 10 rows by removing 2 rows from each block,
# accounting for the upward shift after each deletion.
$ws.Rows("12:13").Delete()
$ws.Rows("22:23").Delete()
$ws.Rows("32:33").Delete()

# Restore the saved selection to match the reverted view state.
$ws.Range("E11").Select()
